$d = $word.ActiveDocument

# --- 1. Split the "Copie : [EDF OA / EDF SEI] ; {#dreal}..." run into several
#        runs so EDF OA / EDF SEI become mustache conditionals -------------
$target = $d.Paragraphs.Last
$full = $target.Range
$start = $full.Start
$end = $full.End - 1
$r = $d.Range($start, $end)

$rpr = '<w:rPr><w:rFonts w:eastAsia="Liberation Sans" w:cs="Liberation Sans" w:ascii="Liberation Sans" w:hAnsi="Liberation Sans"/><w:i/><w:sz w:val="16"/><w:szCs w:val="16"/><w:highlight w:val="cyan"/></w:rPr>'

$run1 = '<w:r>' + $rpr + '<w:t xml:space="preserve">Copie&#160;: </w:t></w:r>'
$run2 = '<w:r>' + $rpr + '<w:t>{#isEDFOA}</w:t></w:r>'
$run3 = '<w:r>' + $rpr + '<w:t>EDF OA</w:t></w:r>'
$run4 = '<w:r>' + $rpr + '<w:t>{/isEDFOA}{#isEDFSEI}</w:t></w:r>'
$run5 = '<w:r>' + $rpr + '<w:t>EDF SEI</w:t></w:r>'
$run6 = '<w:r>' + $rpr + '<w:t>{/isEDFSEI}</w:t></w:r>'
$run7 = '<w:r>' + $rpr + '<w:t>&#160;; {#dreal}DREAL {dreal}{/dreal}{^dreal}DREAL concernée{/dreal} ; CRE</w:t></w:r>'

$body = $run1 + $run2 + $run3 + $run4 + $run5 + $run6 + $run7

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $body + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)

# --- 2. Flip the title-page header image from "behind text" to "in front
#        of text" (behindDoc="1" -> "0") while keeping its square wrap ------
$sec = $d.Sections(1)
$hdr = $sec.Headers(2)
$shp = $hdr.Shapes(1)
$shp.WrapFormat.Type = 3
$shp.WrapFormat.Type = 0
$shp.WrapFormat.Side = 3
